$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.388.24"
$ws.Range("E2").Value = "  -1.67%  "
$ws.Range("D3").Value = "2.300.74"
$ws.Range("E3").Value = "  -3.31%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.71%  "
$ws.Range("E7").Value = "  -2.15%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.611"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0910"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.80%  "
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.968"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.68%  "
$ws.Range("D16").Value = "2.644.55"
$ws.Range("E16").Value = "  -3.28%  "
$ws.Range("D17").Value = "2.298.86"
$ws.Range("E17").Value = "  -3.20%  "
$ws.Range("D18").Value = "42.207.25"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.62%  "
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("B21").Value = "PancakeSwap"
$ws.Range("C21").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "279.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +12.92%  "
$ws.Range("E25").Value = "  -3.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.68%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.21%  "
$ws.Range("E29").Value = "  -2.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "163.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0875"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.31%  "
$ws.Range("E34").Value = "  -5.96%  "
$ws.Range("E35").Value = "  +3.56%  "
$ws.Range("E36").Value = "  -6.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.99%  "
$ws.Range("E38").Value = "  -4.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.38%  "
$ws.Range("E40").Value = "  +2.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.95%  "
$ws.Range("E42").Value = "  -4.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("E44").Value = "  -5.35%  "
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "112.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "77.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.82%  "
$ws.Range("E50").Value = "  -5.38%  "
$ws.Range("D51").Value = "1.616.64"
$ws.Range("E51").Value = "  +2.04%  "
